$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 33.9
    3  = 50.56
    4  = 31.67
    5  = 44.68
    6  = 15.81
    7  = 35.71
    8  = 50.72
    9  = 47.65
    10 = 44.44
    11 = 11.11
    12 = 37.4
    13 = 42.55
    14 = 34.43
    15 = 39.53
    16 = 28.89
    17 = 47.15
    18 = 22.33
    19 = 36.88
    20 = 49.14
    21 = 39.29
    22 = 36.52
    23 = 40.48
    24 = 32.76
    25 = 55
    26 = 30.11
    27 = 32.02
    28 = 29.82
    29 = 40.48
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
